# Auto-generated edit script applying numeric corrections to the
# Leve profit-calculation columns (H,I,J,K,L,M,N) across several sheets,
# per the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1084261.6
$ws.Range("I80").Value = 1422224.5
$ws.Range("K80").Value = 4266673.5
$ws.Range("M80").Value = -4265675.5
$ws.Range("H83").Value = 1084261.6
$ws.Range("I83").Value = 1422224.5
$ws.Range("K83").Value = 12800020.5
$ws.Range("M83").Value = -12795028.5
$ws.Range("H86").Value = 15465638
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 33504800
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 33504800
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -33507046
$ws.Range("H88").Value = 1913.7858
$ws.Range("I88").Value = 1039.4
$ws.Range("K88").Value = 1039.4
$ws.Range("M88").Value = -633.4000000000001
$ws.Range("H89").Value = 15465638
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 33504800
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 167524000
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -167535232
$ws.Range("H91").Value = 1913.7858
$ws.Range("I91").Value = 1039.4
$ws.Range("K91").Value = 1039.4
$ws.Range("M91").Value = 364.5999999999999
$ws.Range("H132").Value = 3656.6487
$ws.Range("I132").Value = 3248.7812
$ws.Range("K132").Value = 9746.3436
$ws.Range("M132").Value = -7216.3436

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2711.0168
$ws.Range("I32").Value = 1609.22
$ws.Range("K32").Value = 1609.22
$ws.Range("M32").Value = -1322.22
$ws.Range("H88").Value = 13890949
$ws.Range("J88").Value = 2601.5
$ws.Range("L88").Value = 2601.5
$ws.Range("N88").Value = -3413.5
$ws.Range("H91").Value = 13890949
$ws.Range("J91").Value = 2601.5
$ws.Range("L91").Value = 2601.5
$ws.Range("N91").Value = -5409.5
$ws.Range("H122").Value = 7754519
$ws.Range("I122").Value = 9526130
$ws.Range("K122").Value = 28578390
$ws.Range("M122").Value = -28575940
$ws.Range("H132").Value = 21278280
$ws.Range("I132").Value = 25642590
$ws.Range("J132").Value = 2268.75
$ws.Range("K132").Value = 76927770
$ws.Range("L132").Value = 6806.25
$ws.Range("M132").Value = -76925240
$ws.Range("N132").Value = -11866.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3955.077
$ws.Range("I86").Value = 3955.077
$ws.Range("K86").Value = 3955.077
$ws.Range("M86").Value = -2832.077
$ws.Range("H89").Value = 3955.077
$ws.Range("I89").Value = 3955.077
$ws.Range("K89").Value = 19775.385
$ws.Range("M89").Value = -14159.385
$ws.Range("H99").Value = 2191.111
$ws.Range("I99").Value = 2191.111
$ws.Range("K99").Value = 2191.111
$ws.Range("M99").Value = -693.1109999999999
$ws.Range("H141").Value = 77800
$ws.Range("J141").Value = 77800
$ws.Range("L141").Value = 77800
$ws.Range("N141").Value = -88160

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3872.6
$ws.Range("I62").Value = 2777
$ws.Range("K62").Value = 2777
$ws.Range("M62").Value = -2153
$ws.Range("H65").Value = 3872.6
$ws.Range("I65").Value = 2777
$ws.Range("K65").Value = 13885
$ws.Range("M65").Value = -10765
$ws.Range("H122").Value = 2051.4285
$ws.Range("J122").Value = 1860.125
$ws.Range("L122").Value = 5580.375
$ws.Range("N122").Value = -10480.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 4672
$ws.Range("I13").Value = 10
$ws.Range("J13").Value = 7003
$ws.Range("K13").Value = 10
$ws.Range("L13").Value = 7003
$ws.Range("M13").Value = 129
$ws.Range("N13").Value = -7281
$ws.Range("H80").Value = 3211.875
$ws.Range("I80").Value = 3568.3333
$ws.Range("K80").Value = 3568.3333
$ws.Range("M80").Value = -2570.3333
$ws.Range("H83").Value = 3211.875
$ws.Range("I83").Value = 3568.3333
$ws.Range("K83").Value = 17841.6665
$ws.Range("M83").Value = -12849.6665
$ws.Range("H99").Value = 6849.1665
$ws.Range("I99").Value = 1274
$ws.Range("K99").Value = 1274
$ws.Range("M99").Value = 972

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4700
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4700
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4700
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4972
$ws.Range("H68").Value = 18666.666
$ws.Range("I68").Value = 13000
$ws.Range("J68").Value = 21500
$ws.Range("K68").Value = 13000
$ws.Range("L68").Value = 21500
$ws.Range("M68").Value = -12251
$ws.Range("N68").Value = -22998
$ws.Range("H71").Value = 18666.666
$ws.Range("I71").Value = 13000
$ws.Range("J71").Value = 21500
$ws.Range("K71").Value = 65000
$ws.Range("L71").Value = 107500
$ws.Range("M71").Value = -61256
$ws.Range("N71").Value = -114988
$ws.Range("H74").Value = 53932
$ws.Range("I74").Value = 50049
$ws.Range("K74").Value = 50049
$ws.Range("M74").Value = -49051
$ws.Range("H77").Value = 53932
$ws.Range("I77").Value = 50049
$ws.Range("K77").Value = 150147
$ws.Range("M77").Value = -145155
$ws.Range("H82").Value = 1237.909
$ws.Range("I82").Value = 920.75
$ws.Range("K82").Value = 920.75
$ws.Range("M82").Value = -559.75
$ws.Range("H85").Value = 1237.909
$ws.Range("I85").Value = 920.75
$ws.Range("K85").Value = 920.75
$ws.Range("M85").Value = 327.25
$ws.Range("H100").Value = 3870.5
$ws.Range("I100").Value = 3827.3333
$ws.Range("K100").Value = 3827.3333
$ws.Range("M100").Value = -3286.3333
$ws.Range("H122").Value = 3250
$ws.Range("I122").Value = 2944.4443
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 8833.332900000001
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -6383.332900000001
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 6138.96
$ws.Range("I132").Value = 3539.1177
$ws.Range("J132").Value = 11663.625
$ws.Range("K132").Value = 10617.3531
$ws.Range("L132").Value = 34990.875
$ws.Range("M132").Value = -8087.3531
$ws.Range("N132").Value = -40050.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 5955478
$ws.Range("J62").Value = 9448.5
$ws.Range("K62").Value = 5955478
$ws.Range("L62").Value = 9448.5
$ws.Range("M62").Value = -5954854
$ws.Range("N62").Value = -10696.5
$ws.Range("I65").Value = 5955478
$ws.Range("J65").Value = 9448.5
$ws.Range("K65").Value = 29777390
$ws.Range("L65").Value = 47242.5
$ws.Range("M65").Value = -29774270
$ws.Range("N65").Value = -53482.5
$ws.Range("H81").Value = 12507084
$ws.Range("I81").Value = 4357
$ws.Range("K81").Value = 8714
$ws.Range("M81").Value = -7653
$ws.Range("H84").Value = 12507084
$ws.Range("I84").Value = 4357
$ws.Range("K84").Value = 43570
$ws.Range("M84").Value = -38266
$ws.Range("H100").Value = 6912.125
$ws.Range("I100").Value = 6912.125
$ws.Range("K100").Value = 13824.25
$ws.Range("M100").Value = -13283.25
$ws.Range("H106").Value = 32599.6
$ws.Range("I106").Value = 24499.5
$ws.Range("K106").Value = 24499.5
$ws.Range("M106").Value = -23237.5
$ws.Range("H109").Value = 31687.5
$ws.Range("J109").Value = 31687.5
$ws.Range("L109").Value = 31687.5
$ws.Range("N109").Value = -34461.5
$ws.Range("H112").Value = 28346.5
$ws.Range("J112").Value = 28346.5
$ws.Range("L112").Value = 28346.5
$ws.Range("N112").Value = -31300.5
$ws.Range("H113").Value = 3777.5386
$ws.Range("I113").Value = 360.8
$ws.Range("K113").Value = 1082.4
$ws.Range("M113").Value = 1087.6
$ws.Range("H126").Value = 1837.2142
$ws.Range("I126").Value = 1683.7273
$ws.Range("K126").Value = 5051.1819
$ws.Range("M126").Value = -2581.1819
$ws.Range("H129").Value = 25000
$ws.Range("J129").Value = 25000
$ws.Range("L129").Value = 25000
$ws.Range("N129").Value = -35000

Write-Output "Applied 207 cell updates."